# "Logica de selecion de fecha dos"
#
# The Recruitment sheet's row 2 had an extra MIDDLE_NAME value ("de Jesus ")
# sitting in column C that doesn't belong there (it shifted every later
# field - LAST_NAME, VACANCY, EMAIL, KEYWORDS, DATE_OF_APPLICATION, NOTE,
# CONSENT_TO_KEEP_DATA, SHORTLIST_NOTE - one column out of alignment with
# their headers). Clearing that stray value re-aligns D2:L2 with their
# header row, and the now-unused "de Jesus " shared string drops out of
# the workbook on save.
#
# The active selection on the Recruitment sheet is also moved from the old
# (now meaningless) I6 over to D6, which is the DATE_OF_APPLICATION-era
# column now lined up under the corrected layout.

$wb = $excel.ActiveWorkbook

$wsRecruitment = $wb.Worksheets.Item("Recruitment")

# Clear the stray middle-name value out of C2 (keeps its style, just drops
# the text/shared-string reference).
$wsRecruitment.Range("C2").Value = ""

# Move the live selection to D6, per the updated date-selection logic.
$wsRecruitment.Range("D6").Select()
